$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Addr,
        [string]$Value
    )
    $rng = $ws.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# --- Price (column D) updates -------------------------------------------------
Set-TextValue "D2"  "249.94"
Set-TextValue "D4"  "5.535"
Set-TextValue "D5"  "0.05643"
Set-TextValue "D6"  "6.452"
Set-TextValue "D8"  "1.036"
Set-TextValue "D10" "0.07311"
Set-TextValue "D11" "0.03123"
Set-TextValue "D12" "0.02916"
Set-TextValue "D13" "0.09265"
Set-TextValue "D14" "0.001662"
Set-TextValue "D17" "0.0005822"
Set-TextValue "D18" "0.006330"
Set-TextValue "D19" "0.005069"
Set-TextValue "D20" "0.001049"
Set-TextValue "D21" "0.0001500"
Set-TextValue "D22" "3.977"
Set-TextValue "D23" "3.380"
Set-TextValue "D27" "0.0003059"
Set-TextValue "D40" "0.04146"
Set-TextValue "D41" "0.006880"
Set-TextValue "D42" "0.003500"
Set-TextValue "D43" "0.1043"
Set-TextValue "D44" "0.009377"
Set-TextValue "D45" "0.00005651"
Set-TextValue "D47" "0.6802"
Set-TextValue "D48" "0.01579"

# --- Row 17 (One / ONE) volume label now flags "Worst in 24h" -----------------
$ws.Range("E17").Value = "16OneONEWorstin24h"

# --- Rows 41 & 43: BKEXToken and KickToken swap positions ----------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Write-Output "edits applied"
